$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.49%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.674"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.57%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.62%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.033"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.753"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.36%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.525"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.42%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.963"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.72%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9209"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.62%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1259"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.07%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1945"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.29%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.307"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-7.97%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09296"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.68%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03700"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "7.34%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1052"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.21%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001300"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.17%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006248"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.25%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.361"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.05%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3474"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.49%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.41%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "10.05%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04442"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.24%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001262"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.06%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004291"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.60%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "13.84%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02860"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "15.97%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05470"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.04%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007784"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.06%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009951"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "12.99%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1417"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.24%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002230"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.71%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01182"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "13.57%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006771"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.90%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.03%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003019"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-13.88%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002282"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "33.87%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.03%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.03%"
